# Add a new row of data ("Part3" / "The Keys to Success") below the existing
# table, style the new description cell, and highlight it with conditional
# formatting (duplicate-values / "Keys to Success" style), matching the
# "add practice and tables" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content -------------------------------------------------
# Set B3 before A3 so the shared-string table ends up in the same order
# ("The Keys to Success" then "Part3") as the target workbook.
$ws.Range("B3").Value = "The Keys to Success"
$ws.Range("A3").Value = "Part3"

# --- Formatting for the new description cell (B3) ---------------------
$ws.Range("B3").Font.Size = 12
$ws.Range("B3").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 15.5

# --- Selection, as left by the editing session -------------------------
$ws.Range("A4").Select()

# --- Conditional formatting: 3 "duplicate values" rules on B3 ----------
# (mirrors applying the built-in "Duplicate Values" highlight a few times)
$fc1 = $ws.Range("B3").FormatConditions.AddUniqueValues()
$fc1.DupeUnique = 1
$fc2 = $ws.Range("B3").FormatConditions.AddUniqueValues()
$fc2.DupeUnique = 1
$fc3 = $ws.Range("B3").FormatConditions.AddUniqueValues()
$fc3.DupeUnique = 1

# Apply the "light red fill with dark red text" look to each rule. The
# dxf records get created in the order the formats are applied below, so
# style fc3 first, then fc2, then fc1 - this lines up the resulting
# dxfId values (2,1,0) with priorities (3,2,1) the same way as the
# target file.
$fc3.Font.Color = 393372
$fc3.Interior.Color = 13551615
$fc2.Font.Color = 393372
$fc2.Interior.Color = 13551615
$fc1.Font.Color = 393372
$fc1.Interior.Color = 13551615

$fc1.Priority = 3
$fc2.Priority = 2
$fc3.Priority = 1
